$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 20019
$ws.Range("I21").Value = 20019
$ws.Range("K21").Value = 20019
$ws.Range("M21").Value = -19551
# Row 23
$ws.Range("H23").Value = 20019
$ws.Range("I23").Value = 20019
$ws.Range("K23").Value = 20019
$ws.Range("M23").Value = -19785
# Row 40
$ws.Range("H40").Value = 6500.4287
$ws.Range("I40").Value = 6375
$ws.Range("J40").Value = 6667.6665
$ws.Range("K40").Value = 6375
$ws.Range("L40").Value = 6667.6665
$ws.Range("M40").Value = -6200
$ws.Range("N40").Value = -7017.6665
# Row 41
$ws.Range("H41").Value = 434.5
$ws.Range("I41").Value = 391.16666
$ws.Range("K41").Value = 391.16666
$ws.Range("M41").Value = 48.83334000000002
# Row 76
$ws.Range("H76").Value = 3495.6667
$ws.Range("I76").Value = 3498.5
$ws.Range("J76").Value = 3490
$ws.Range("K76").Value = 3498.5
$ws.Range("L76").Value = 3490
$ws.Range("M76").Value = -3183.5
$ws.Range("N76").Value = -4120
# Row 79
$ws.Range("H79").Value = 3495.6667
$ws.Range("I79").Value = 3498.5
$ws.Range("J79").Value = 3490
$ws.Range("K79").Value = 3498.5
$ws.Range("L79").Value = 3490
$ws.Range("M79").Value = -2406.5
$ws.Range("N79").Value = -5674
# Row 80
$ws.Range("H80").Value = 593.53845
$ws.Range("I80").Value = 203.71428
$ws.Range("K80").Value = 611.14284
$ws.Range("M80").Value = 386.85716
# Row 83
$ws.Range("H83").Value = 593.53845
$ws.Range("I83").Value = 203.71428
$ws.Range("K83").Value = 1833.42852
$ws.Range("M83").Value = 3158.57148
# Row 133
$ws.Range("H133").Value = 108978
$ws.Range("J133").Value = 108978
$ws.Range("L133").Value = 108978
$ws.Range("N133").Value = -119098
# Row 135
$ws.Range("H135").Value = 749.53845
$ws.Range("I135").Value = 395.41666
$ws.Range("J135").Value = 4999
$ws.Range("K135").Value = 3558.74994
$ws.Range("L135").Value = 44991
$ws.Range("M135").Value = -1023.74994
$ws.Range("N135").Value = -50061
# Row 137
$ws.Range("H137").Value = 2730.3845
$ws.Range("I137").Value = 3484.8333
$ws.Range("K137").Value = 10454.4999
$ws.Range("M137").Value = -7904.499899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2849.1875
$ws.Range("I61").Value = 1955.7142
$ws.Range("J61").Value = 3544.111
$ws.Range("K61").Value = 1955.7142
$ws.Range("L61").Value = 3544.111
$ws.Range("M61").Value = -1743.7142
$ws.Range("N61").Value = -3968.111
# Row 132
$ws.Range("H132").Value = 3448.8333
$ws.Range("I132").Value = 3564.3333
$ws.Range("K132").Value = 10692.9999
$ws.Range("M132").Value = -8162.999899999999
# Row 136
$ws.Range("H136").Value = 2849.1875
$ws.Range("I136").Value = 1955.7142
$ws.Range("J136").Value = 3544.111
$ws.Range("K136").Value = 5867.142599999999
$ws.Range("L136").Value = 10632.333
$ws.Range("M136").Value = -3317.142599999999
$ws.Range("N136").Value = -15732.333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 847.625
$ws.Range("I22").Value = 847.625
$ws.Range("K22").Value = 847.625
$ws.Range("M22").Value = -674.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 875
$ws.Range("I22").Value = 833.3333
$ws.Range("K22").Value = 833.3333
$ws.Range("M22").Value = -483.3333
# Row 31
$ws.Range("H31").Value = 4170125
$ws.Range("I31").Value = 3108.6428
$ws.Range("J31").Value = 7816264
$ws.Range("K31").Value = 3108.6428
$ws.Range("L31").Value = 7816264
$ws.Range("M31").Value = -2813.6428
$ws.Range("N31").Value = -7816854
# Row 34
$ws.Range("H34").Value = 4170125
$ws.Range("I34").Value = 3108.6428
$ws.Range("J34").Value = 7816264
$ws.Range("K34").Value = 3108.6428
$ws.Range("L34").Value = 7816264
$ws.Range("M34").Value = -2906.6428
$ws.Range("N34").Value = -7816668
# Row 58
$ws.Range("H58").Value = 1785.28
$ws.Range("I58").Value = 1224.7858
$ws.Range("K58").Value = 1224.7858
$ws.Range("M58").Value = -1021.7858
# Row 105
$ws.Range("H105").Value = 2643.25
$ws.Range("I105").Value = 2672
$ws.Range("J105").Value = 2499.5
$ws.Range("K105").Value = 2672
$ws.Range("L105").Value = 2499.5
$ws.Range("M105").Value = -925
$ws.Range("N105").Value = -5993.5
# Row 134
$ws.Range("H134").Value = 4603.9565
$ws.Range("I134").Value = 4563.227
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 13689.681
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -11154.681
$ws.Range("N134").Value = -21570
# Row 136
$ws.Range("H136").Value = 1785.28
$ws.Range("I136").Value = 1224.7858
$ws.Range("K136").Value = 3674.3574
$ws.Range("M136").Value = -1124.3574

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 2666.6667
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
# Row 21
$ws.Range("H21").Value = 2225.1428
$ws.Range("I21").Value = 50
$ws.Range("J21").Value = 3095.2
$ws.Range("K21").Value = 150
$ws.Range("L21").Value = 9285.599999999999
$ws.Range("M21").Value = 23
$ws.Range("N21").Value = -9631.599999999999
# Row 86
$ws.Range("H86").Value = 1946.4286
$ws.Range("J86").Value = 2585
$ws.Range("L86").Value = 7755
$ws.Range("N86").Value = -10127
# Row 89
$ws.Range("H89").Value = 1946.4286
$ws.Range("J89").Value = 2585
$ws.Range("L89").Value = 23265
$ws.Range("N89").Value = -35121
# Row 137
$ws.Range("H137").Value = 1586.5
$ws.Range("I137").Value = 1781
$ws.Range("J137").Value = 1197.5
$ws.Range("K137").Value = 5343
$ws.Range("L137").Value = 3592.5
$ws.Range("M137").Value = -243
$ws.Range("N137").Value = -13792.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 9999.666999999999
$ws.Range("I70").Value = 9999
$ws.Range("K70").Value = 9999
$ws.Range("M70").Value = -9729
# Row 73
$ws.Range("H73").Value = 9999.666999999999
$ws.Range("I73").Value = 9999
$ws.Range("K73").Value = 9999
$ws.Range("M73").Value = -9063

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 34000016
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 51000000
$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 51000000
$ws.Range("M2").Value = 62
$ws.Range("N2").Value = -51000224
# Row 7
$ws.Range("H7").Value = 1910.4445
$ws.Range("I7").Value = 1281.1818
$ws.Range("K7").Value = 1281.1818
$ws.Range("M7").Value = -1169.1818
# Row 40
$ws.Range("H40").Value = 16702.285
$ws.Range("I40").Value = 19729.363
$ws.Range("J40").Value = 5603
$ws.Range("K40").Value = 19729.363
$ws.Range("L40").Value = 5603
$ws.Range("M40").Value = -19593.363
$ws.Range("N40").Value = -5875
# Row 46
$ws.Range("H46").Value = 2174.0715
$ws.Range("I46").Value = 1834.3
$ws.Range("J46").Value = 3023.5
$ws.Range("K46").Value = 1834.3
$ws.Range("L46").Value = 3023.5
$ws.Range("M46").Value = -1646.3
$ws.Range("N46").Value = -3399.5
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 107
$ws.Range("H107").Value = 4996
$ws.Range("I107").Value = 4996
$ws.Range("K107").Value = 4996
$ws.Range("M107").Value = -3076
# Row 126
$ws.Range("H126").Value = 1910.4445
$ws.Range("I126").Value = 1281.1818
$ws.Range("K126").Value = 3843.5454
$ws.Range("M126").Value = -1373.5454

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 106
$ws.Range("H106").Value = 39999.5
$ws.Range("J106").Value = 39999.5
$ws.Range("L106").Value = 39999.5
$ws.Range("N106").Value = -42523.5
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 122
$ws.Range("H122").Value = 27780792
$ws.Range("I122").Value = 3629.4
$ws.Range("K122").Value = 10888.2
$ws.Range("M122").Value = -8438.200000000001
# Row 136
$ws.Range("H136").Value = 3821.5
$ws.Range("I136").Value = 4045.6365
$ws.Range("K136").Value = 12136.9095
$ws.Range("M136").Value = -9586.9095
